# Update the "想去人数" (interest count, column F) figures in the
# 江西-漫展信息 workbook to the freshly scraped values.
# These numbers were refreshed on both the "展览" sheet and the
# "全部类型" aggregate sheet (which mirrors the same rows); the
# "演出" and "本地生活" sheets are untouched since they contain no data.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 2095
    "F3"  = 134
    "F6"  = 1757
    "F8"  = 738
    "F17" = 139
    "F18" = 4081
    "F21" = 453
    "F23" = 956
    "F24" = 1022
    "F28" = 1851
    "F29" = 49
    "F32" = 184
    "F33" = 20
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
